$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row before row 8; everything from row 8 down shifts to row 9 down.
$ws.Rows(8).Insert()

# 2. Copy the formatting (borders/number format/alignment) from row 7 into the
#    freshly inserted row 8 so it matches the look of the table above it.
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Populate the new row 8 with the new experiment entry.
$ws.Range("A8").Value = "CodeGPTPy"
$ws.Range("B8").Value = 128
$ws.Range("C8").Value = "Early Stopping"
$ws.Range("D8").Value = "Adam"
$ws.Range("E8").Value = 0
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = "run 5 times each"

# 4. The row that used to be row 11 (now row 12) had an in-progress entry;
#    clear its leading columns since that run was dropped.
$ws.Range("A12:E12").ClearContents()

# 5. Restore the active selection.
$ws.Range("I9").Select() | Out-Null
